$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2445987" "1122334"
Replace-Text "01" "07"
Replace-Text "мая" "марта"
Replace-Text '"15 ВАХАУ МАРУ"' '"СИНЕГОРСК"'
Replace-Text "172316" "021026"
Replace-Text "--" "9261061"
Replace-Text "Внеочередное освидетельствование в связи со сменой судовладельца" "Первоначальное освидетельствование"
Replace-Text "01.05.2024" "07.03.2024"
Replace-Text "генерального директора Котлярчука О. Е." "капитана Бахтина Ю. Г."
Replace-Text "Устава" "Кодекса торгового мореплавания (КТМ РФ)"
Replace-Text "Бахтин Ю. Г., 89611823023, Bavenit.master@rusgeology.ru" "Котлярчук О. Е., +79520528053, kotlyarchuk@gmail.com"
Replace-Text "Архангельск, Архангельская область" "Мурманск, Мурманская область"
Replace-Text "23.01.2023" "11.10.2023"
Replace-Text "22.01.2022 без замечаний" "21.11.2023 несоответствий нет"
Replace-Text "О. Е. Котлярчук" "Ю. Г. Бахтин"
